# Auto-generated Excel COM-interop script
# Converts column B start-date strings from dot-separated (2024.03.02)
# to dash-separated (2024-03-02) format across all three data sheets,
# and bumps a handful of '想去人数' (interest count) values in column F
# to match the output snapshot regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")

# Column B: '.' -> '-' date separator, rows 2..43
$dates = @{
    2 = "2024-03-02"
    3 = "2024-03-02"
    4 = "2024-03-03"
    5 = "2024-03-09"
    6 = "2024-03-09"
    7 = "2024-03-09"
    8 = "2024-03-16"
    9 = "2024-03-16"
    10 = "2024-03-16"
    11 = "2024-03-16"
    12 = "2024-03-16"
    13 = "2024-03-16"
    14 = "2024-03-16"
    15 = "2024-03-16"
    16 = "2024-03-16"
    17 = "2024-03-16"
    18 = "2024-03-17"
    19 = "2024-03-23"
    20 = "2024-03-23"
    21 = "2024-03-23"
    22 = "2024-03-23"
    23 = "2024-03-24"
    24 = "2024-03-24"
    25 = "2024-03-24"
    26 = "2024-03-30"
    27 = "2024-03-30"
    28 = "2024-03-30"
    29 = "2024-03-30"
    30 = "2024-04-04"
    31 = "2024-04-04"
    32 = "2024-04-04"
    33 = "2024-04-05"
    34 = "2024-04-05"
    35 = "2024-04-05"
    36 = "2024-04-13"
    37 = "2024-04-13"
    38 = "2024-04-20"
    39 = "2024-04-20"
    40 = "2024-04-20"
    41 = "2024-04-20"
    42 = "2024-07-20"
    43 = "2024-07-20"
}

$rng = $ws.Range("B2:B43")
$rng.NumberFormat = "@"   # force text so Excel does not reparse these as dates
foreach ($row in $dates.Keys) {
    $ws.Range("B$row").Value = $dates[$row]
}
$rng.ClearFormats()       # restore default (General) styling, no visual change since cells are text

# Column F: updated '想去人数' counts
$ws.Range("F6").Value = 2941
$ws.Range("F8").Value = 1949
$ws.Range("F19").Value = 7050
$ws.Range("F21").Value = 1752
$ws.Range("F29").Value = 938

$ws = $wb.Worksheets.Item("演出")

# Column B: '.' -> '-' date separator, rows 2..9
$dates = @{
    2 = "2024-03-03"
    3 = "2024-03-08"
    4 = "2024-03-08"
    5 = "2024-03-16"
    6 = "2024-03-23"
    7 = "2024-04-21"
    8 = "2024-04-27"
    9 = "2024-05-01"
}

$rng = $ws.Range("B2:B9")
$rng.NumberFormat = "@"   # force text so Excel does not reparse these as dates
foreach ($row in $dates.Keys) {
    $ws.Range("B$row").Value = $dates[$row]
}
$rng.ClearFormats()       # restore default (General) styling, no visual change since cells are text

# Column F: updated '想去人数' counts
$ws.Range("F6").Value = 8

$ws = $wb.Worksheets.Item("全部类型")

# Column B: '.' -> '-' date separator, rows 2..49
$dates = @{
    2 = "2024-03-02"
    3 = "2024-03-02"
    4 = "2024-03-03"
    5 = "2024-03-03"
    6 = "2024-03-08"
    7 = "2024-03-08"
    8 = "2024-03-09"
    9 = "2024-03-09"
    10 = "2024-03-09"
    11 = "2024-03-16"
    12 = "2024-03-16"
    13 = "2024-03-16"
    14 = "2024-03-16"
    15 = "2024-03-16"
    16 = "2024-03-16"
    17 = "2024-03-16"
    18 = "2024-03-16"
    19 = "2024-03-16"
    20 = "2024-03-16"
    21 = "2024-03-17"
    22 = "2024-03-23"
    23 = "2024-03-23"
    24 = "2024-03-23"
    25 = "2024-03-23"
    26 = "2024-03-23"
    27 = "2024-03-24"
    28 = "2024-03-24"
    29 = "2024-03-24"
    30 = "2024-03-30"
    31 = "2024-03-30"
    32 = "2024-03-30"
    33 = "2024-03-30"
    34 = "2024-04-04"
    35 = "2024-04-04"
    36 = "2024-04-05"
    37 = "2024-04-05"
    38 = "2024-04-05"
    39 = "2024-04-13"
    40 = "2024-04-13"
    41 = "2024-04-20"
    42 = "2024-04-20"
    43 = "2024-04-20"
    44 = "2024-04-20"
    45 = "2024-04-21"
    46 = "2024-04-27"
    47 = "2024-05-01"
    48 = "2024-07-20"
    49 = "2024-07-20"
}

$rng = $ws.Range("B2:B49")
$rng.NumberFormat = "@"   # force text so Excel does not reparse these as dates
foreach ($row in $dates.Keys) {
    $ws.Range("B$row").Value = $dates[$row]
}
$rng.ClearFormats()       # restore default (General) styling, no visual change since cells are text

# Column F: updated '想去人数' counts
$ws.Range("F9").Value = 2941
$ws.Range("F11").Value = 1949
$ws.Range("F22").Value = 7050
$ws.Range("F24").Value = 1752
$ws.Range("F25").Value = 8
$ws.Range("F33").Value = 938

